$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append additional FOMC cycle date rows (69-92), continuing the
# --- existing "FOMC Date Start" / "FOMC Date End" table in columns A/B.
$ws.Range("A69").Value = 40162
$ws.Range("B69").Value = 40163
$ws.Range("A70").Value = 40120
$ws.Range("B70").Value = 40121
$ws.Range("A71").Value = 40078
$ws.Range("B71").Value = 40079
$ws.Range("A72").Value = 40036
$ws.Range("B72").Value = 40037
$ws.Range("A73").Value = 39987
$ws.Range("B73").Value = 39988
$ws.Range("A74").Value = 39931
$ws.Range("B74").Value = 39932
$ws.Range("A75").Value = 39889
$ws.Range("B75").Value = 39890
$ws.Range("A76").Value = 39840
$ws.Range("B76").Value = 39841
$ws.Range("A77").Value = 39797
$ws.Range("B77").Value = 39798
$ws.Range("A78").Value = 39749
$ws.Range("B78").Value = 39750
$ws.Range("A79").Value = 39707
$ws.Range("A80").Value = 39665
$ws.Range("A81").Value = 39623
$ws.Range("B81").Value = 39624
$ws.Range("A82").Value = 39567
$ws.Range("B82").Value = 39568
$ws.Range("A83").Value = 39525
$ws.Range("A84").Value = 39476
$ws.Range("B84").Value = 39477
$ws.Range("A85").Value = 39427
$ws.Range("A86").Value = 39385
$ws.Range("B86").Value = 39386
$ws.Range("A87").Value = 39343
$ws.Range("A88").Value = 39301
$ws.Range("A89").Value = 39260
$ws.Range("B89").Value = 39261
$ws.Range("A90").Value = 39211
$ws.Range("A91").Value = 39161
$ws.Range("B91").Value = 39162
$ws.Range("A92").Value = 39112
$ws.Range("B92").Value = 39113

# Make sure the newly-added cells (including the intentionally blank
# B cells) carry the same date formatting as the rest of the table.
$ws.Range("A69:B92").NumberFormat = $ws.Range("A68").NumberFormat

# --- Update the saved view state of the sheet (scroll position, zoom
# --- level and selected cell) to match the new working position.
$ws.Range("A23").Select()
$excel.ActiveWindow.Zoom = 65
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

# Standard column width drifts very slightly when the zoom level
# changes in real Excel; reflect that here as well.
$ws.StandardWidth = 10.4296875
